# Assign Shift.xlsx - add a "ShiftName" dropdown (data validation list) fed
# from a new lookup sheet containing the list of available shifts.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Add the lookup sheet (Sheet2) with the shift names -------------------
# Add it directly after Sheet1 (Worksheet.Move isn't supported by the host,
# so positioning has to happen at creation time via Before/After).
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$shiftNames = @(
    "Early Shift(06:00-14:30)",
    "Saturday Shift(06:00-11:00)",
    "Morning Shift(07:00-15:30)",
    "General Shift(08:00-16:30)",
    "Sat Shift(08:00-13:00)",
    "Day Shift(09:00-17:30)",
    "Extended Day Shift(10:00-18:30)",
    "Mid Shift(12:00-20:30)",
    "Afternoon Shift(14:00-22:30)",
    "Evening Shift(15:00-00:00)",
    "Late Evening Shift(16:00-01:00)",
    "Late Shift(17:00-02:00)",
    "Night Shift(18:00-03:00)",
    "Extended Night Shift(18:30-03:30)",
    "Late Night Shift(19:30-04:30)",
    "Overnight Shift(21:00-06:00)",
    "Midnight Shift(22:00-07:00)",
    "Weekly Off(00:00-00:00)"
)

for ($i = 0; $i -lt $shiftNames.Count; $i++) {
    $ws2.Cells.Item($i + 1, 1).Value = $shiftNames[$i]
}

$ws2.Columns.Item(1).ColumnWidth = 30.125

# --- Back on Sheet1: formatting + the new data row -------------------------
$ws1.Columns.Item(1).ColumnWidth = 10.125
$ws1.Columns.Item(3).ColumnWidth = 14.25

# Row 2 holds the next entry to fill in; format the Date cell as a (built-in)
# short date, numFmtId 14.
$ws1.Range("A2").NumberFormat = "mm-dd-yy"

# Data validation dropdown for the ShiftName column, sourced from Sheet2.
$range = $ws1.Range("C1:C1048576")
$range.Validation.Delete()
$range.Validation.Add(3, 1, 1, "=Sheet2!`$A`$1:`$A`$18")
$range.Validation.IgnoreBlank = $true
$range.Validation.InCellDropdown = $true
$range.Validation.ShowInput = $true
$range.Validation.ShowError = $true

$ws1.Select()
$ws1.Range("E6").Select()
